$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- 1. TestingResults: clear the autofilter criteria so all rows show again ---
$wsResults = $wb.Worksheets.Item("TestingResults")
$wsResults.ShowAllData()
$wsResults.Activate()
$wsResults.Range("F17").Select()

# --- 2. AbnormalGuidelines: keep selection, just stop being the active tab (handled by later activation) ---
$wsAbnormal = $wb.Worksheets.Item("AbnormalGuidelines")
$wsAbnormal.Activate()
$wsAbnormal.Range("I7").Select()

# --- 3. ReportTable: replace A:B list with the compact D:E "abnormal" list, drop D:E ---
$wsReport = $wb.Worksheets.Item("ReportTable")
$wsReport.Range("A1:B101").ClearContents()
$wsReport.Range("D1:E57").Cut($wsReport.Range("A1"))

# --- 4. Delete the now-redundant GuideLine5 sheet ---
$wsGuide5 = $wb.Worksheets.Item("GuideLine5")
$wsGuide5.Delete()

# --- 5. ReportTable becomes the active tab/selection ---
$wsReport.Activate()
$wsReport.Range("F10").Select()
